$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.199.14"
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = "'2.487.24"
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'519.40"
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = "'131.89"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").Value = "'2.521.91"
$ws.Range("E9").Value = '  +2.37%  '
$ws.Range("D10").Value = "'0.0973"
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = "'5.20"
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").Value = "'0.332"
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("D14").Value = "'2.964.01"
$ws.Range("E14").Value = '  +2.10%  '
$ws.Range("D15").Value = "'58.130.96"
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = "'22.11"
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").Value = "'2.513.69"
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = "'10.71"
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'321.26"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = "'5.99"
$ws.Range("E22").Value = '  +4.72%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = "'64.28"
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").Value = "'0.401"
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("D28").Value = "'7.34"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").Value = "'0.0₃0748"
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("D30").Value = "'168.02"
$ws.Range("E30").Value = '  +1.85%  '
$ws.Range("D31").Value = "'1.71"
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").Value = "'6.24"
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = "'0.994"
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").Value = "'1.26"
$ws.Range("E37").Value = '  -2.69%  '
$ws.Range("D38").Value = "'3.93"
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("D39").Value = "'36.95"
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").Value = "'0.768"
$ws.Range("E41").Value = '  -2.04%  '
$ws.Range("D42").Value = "'276.88"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'5.09"
$ws.Range("E43").Value = '  +2.46%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = "'3.44"
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("D47").Value = "'121.54"
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").Value = "'0.0499"
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("D49").Value = "'17.79"
$ws.Range("E49").Value = '  +0.69%  '
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").Value = "'16.91"
$ws.Range("E51").Value = '  -0.06%  '
